# Add a new "Italy" worksheet (test data), copied from the existing
# "Slovakia" sheet, with the market name / user-story values swapped in.

$wb = $excel.ActiveWorkbook

$source = $wb.Worksheets.Item("Slovakia")

# Selecting the full sheet on the source mirrors what Excel leaves behind
# on the sheet you copied *from* once the new tab becomes active.
$source.Cells.Select() | Out-Null

# Copy "Slovakia" to the end of the tab strip (after the last sheet).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$source.Copy($null, $lastSheet)

$italy = $wb.Worksheets.Item($wb.Worksheets.Count)
$italy.Name = "Italy"

# Update the market / user-story cells for the new sheet.
$italy.Range("B2").Value = "Italy Market"
$italy.Range("B4").Value = "NGC-3145/T2219"

# Make the new sheet the active one, with B4 selected.
$italy.Activate() | Out-Null
$italy.Range("B4").Select() | Out-Null
